# Apply the "new version with timestamp" update to the DaySale / missing-items
# report: six new product rows are inserted (keeping the existing
# alphabetical ordering), the grand total is refreshed, and the
# generated-at timestamp in the footer is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the six new rows. We work from the bottom of the table
#    upwards so that row numbers we still need to reference below stay
#    valid while we work. Each insert copies the formatting (styles +
#    merged cells) of the row immediately above the insertion point so
#    the new row matches the rest of the table.
# ---------------------------------------------------------------------

function Insert-RowLike([int]$anchorRow, [int]$templateRow) {
    $ws.Rows("$($templateRow):$($templateRow)").Copy()
    $ws.Rows("$($anchorRow):$($anchorRow)").Insert()
}

# before original row 21 (ماكينه حلاقه جليت فليكتور) -> "كالونا"
Insert-RowLike 21 20
# before the row now holding "كالونا" -> "سرنجات 5 سم"
Insert-RowLike 21 20
# before original row 20 (سرنجات 3 سم) -> "جهاز محلول"
Insert-RowLike 20 19
# before original row 17 (ZURCAL ...) -> "VISCERALGINE ..."
Insert-RowLike 17 16
# before original row 16 (TRILLERG ...) -> "STREPTOQUIN ..."
Insert-RowLike 16 15
# before original row 8 (ANTINAL ...) -> "AMRIZOLE ..."
Insert-RowLike 8 7

# ---------------------------------------------------------------------
# 2. Helper to write one data row. L and P are backed by numeric-looking
#    number formats (#,##0.##;... and 0.00) in this sheet, but the
#    source report always stores these particular columns as plain text
#    ("1", "47.0000", ...). Flip the format to Text, write the value,
#    then restore the original format code so the cell keeps rendering
#    with the same style it had before while remaining text internally.
# ---------------------------------------------------------------------

function Set-TextValue($range, [string]$value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

function Set-DataRow([int]$r, [string]$name, [string]$balance, [string]$reorder, [string]$price, [string]$sellPrice, [string]$txCount) {
    $ws.Range("C$r").Value = $name
    $ws.Range("H$r").Value = $balance
    Set-TextValue $ws.Range("L$r") $reorder
    $ws.Range("N$r").Value = $price
    Set-TextValue $ws.Range("P$r") $sellPrice
    $ws.Range("Q$r").Value = $txCount
}

# ---------------------------------------------------------------------
# 3. Re-write every data row (7-27) so the final content, including the
#    serial numbers in column A, matches the post-edit report exactly.
# ---------------------------------------------------------------------

$ws.Range("A7").Value = 1
Set-DataRow 7 "ALPHINTERN 30 F.C.TABS" "0:2" "1" "87.00" "28.7100" "0:1"

$ws.Range("A8").Value = 2
Set-DataRow 8 "AMRIZOLE 500MG VIAL 100 ML" "7:0" "1" "47.00" "47.0000" "1:0"

$ws.Range("A9").Value = 3
Set-DataRow 9 "ANTINAL 220MG/5ML 60ML SUSP." "2:0" "1" "24.00" "24.0000" "1:0"

$ws.Range("A10").Value = 4
Set-DataRow 10 "CATAFLAM 75MG/3ML 6 AMP." "0:4" "1" "120.00" "19.2000" "0:1"

$ws.Range("A11").Value = 5
Set-DataRow 11 "CEVA-FRESH TAB" "4:0" "0" "12.00" "12.0000" "1:0"

$ws.Range("A12").Value = 6
Set-DataRow 12 "CONVENTIN XR 600MG 30 TABS." "0:0" "1" "198.00" "65.3400" "0:1"

$ws.Range("A13").Value = 7
Set-DataRow 13 "DENSITIN 30 CAPS" "0:0" "1" "96.00" "96.0000" "1:0"

$ws.Range("A14").Value = 8
Set-DataRow 14 "EXTRAUMA DNA FORTE TOPICAL GEL 25 GM" "1:0" "1" "41.00" "41.0000" "1:0"

$ws.Range("A15").Value = 9
Set-DataRow 15 "FUSI-ZON CREAM 15 GM" "1:0" "1" "48.00" "96.0000" "2:0"

$ws.Range("A16").Value = 10
Set-DataRow 16 "RIVO 320MG 20*10 TABS" "0:8" "1" "141.00" "14.1000" "0:2"

$ws.Range("A17").Value = 11
Set-DataRow 17 "STREPTOQUIN 20 TABLETS" "3:0" "1" "46.00" "23.0000" "0:1"

$ws.Range("A18").Value = 12
Set-DataRow 18 "TRILLERG EYE DROPS 10 ML" "1:0" "1" "24.00" "24.0000" "1:0"

$ws.Range("A19").Value = 13
Set-DataRow 19 "VISCERALGINE 5MG/2ML IM IV 6 AMPOULES" "0:5" "1" "90.00" "14.4000" "0:1"

$ws.Range("A20").Value = 14
Set-DataRow 20 "ZURCAL 40MG 14 GASTRO RESISTANT TAB" "3:0" "1" "96.00" "96.0000" "1:0"

$ws.Range("A21").Value = 15
Set-DataRow 21 "اختبار حمل بيبي تشك " "14:0" "0" "25.00" "25.0000" "1:0"

$ws.Range("A22").Value = 16
Set-DataRow 22 "جنتيانا " "5:0" "0" "15.00" "15.0000" "1:0"

$ws.Range("A23").Value = 17
Set-DataRow 23 "جهاز محلول " "64:0" "0" "20.00" "20.0000" "1:0"

$ws.Range("A24").Value = 18
Set-DataRow 24 "سرنجات 3 سم" "0:0" "0" "2.00" "2.0000" "1:0"

$ws.Range("A25").Value = 19
Set-DataRow 25 "سرنجات 5 سم" "0:0" "0" "3.00" "3.0000" "1:0"

$ws.Range("A26").Value = 20
Set-DataRow 26 "كالونا " "0:0" "0" "15.00" "15.0000" "1:0"

$ws.Range("A27").Value = 21
Set-DataRow 27 "ماكينه حلاقه جليت فليكتور" "21:0" "0" "15.00" "15.0000" "1:0"

# ---------------------------------------------------------------------
# 4. Refresh the grand total and the generated-at footer.
# ---------------------------------------------------------------------

$ws.Range("P28").Value = 695.75
$ws.Range("A29").Value = "Tuesday, 5 August, 2025 11:05 AM"

# ---------------------------------------------------------------------
# 5. Match the original report's row heights exactly (rows keep
#    alternating 25.5/24.75 heights that come from the source tool, not
#    from any simple odd/even rule, so set them explicitly).
# ---------------------------------------------------------------------

$heights = @{
    7 = 25.5; 8 = 24.75; 9 = 25.5; 10 = 24.75; 11 = 25.5; 12 = 25.5; 13 = 24.75;
    14 = 25.5; 15 = 24.75; 16 = 25.5; 17 = 25.5; 18 = 24.75; 19 = 25.5; 20 = 24.75;
    21 = 25.5; 22 = 25.5; 23 = 24.75; 24 = 25.5; 25 = 24.75; 26 = 25.5; 27 = 25.5;
    28 = 24.75; 29 = 16.5
}
foreach ($r in $heights.Keys) {
    $ws.Rows("$($r):$($r)").RowHeight = $heights[$r]
}

Write-Output "done"
